$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text must be forced to remain text
# (Price column values like "207.98" would otherwise be auto-converted to numbers)
$textCells = @("D5", "D6", "D8", "D9", "D11", "D14", "D15", "D19", "D20", "D23", "D25", "D27", "D28", "D32", "D34", "D35", "D39", "D40", "D44", "D45", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '27.384.61'
$ws.Range('E2').Value = '  -1.08%  '
$ws.Range('D3').Value = '1.564.35'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '207.98'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('D6').Value = '0.498'
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '21.83'
$ws.Range('E8').Value = '  -2.28%  '
$ws.Range('D9').Value = '0.248'
$ws.Range('E9').Value = '  -2.22%  '
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('D11').Value = '0.0867'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').Value = '1.786.47'
$ws.Range('E12').Value = '  -1.19%  '
$ws.Range('D13').Value = '1.575.32'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').Value = '3.82'
$ws.Range('E14').Value = '  -1.25%  '
$ws.Range('D15').Value = '0.516'
$ws.Range('E15').Value = '  -2.92%  '
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = '27.399.86'
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.0₃0688'
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '211.71'
$ws.Range('E19').Value = '  -2.68%  '
$ws.Range('D20').Value = '7.26'
$ws.Range('E20').Value = '  -1.20%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('E22').Value = '  -1.29%  '
$ws.Range('D23').Value = '9.52'
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('D25').Value = '153.34'
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = '6.71'
$ws.Range('E27').Value = '  +0.18%  '
$ws.Range('D28').Value = '14.98'
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('E29').Value = '  -2.08%  '
$ws.Range('E30').Value = '  -0.27%  '
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').Value = '3.20'
$ws.Range('E32').Value = '  -1.23%  '
$ws.Range('D33').Value = '1.363.33'
$ws.Range('E33').Value = '  -1.04%  '
$ws.Range('D34').Value = '2.95'
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('D35').Value = '1.54'
$ws.Range('E35').Value = '  +1.31%  '
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('E38').Value = '  +0.69%  '
$ws.Range('D39').Value = '0.531'
$ws.Range('E39').Value = '  -1.37%  '
$ws.Range('D40').Value = '0.821'
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  -0.46%  '
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').Value = '63.98'
$ws.Range('E44').Value = '  +0.48%  '
$ws.Range('D45').Value = '5.28'
$ws.Range('E45').Value = '  +0.92%  '
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('D47').Value = '1.699.53'
$ws.Range('E47').Value = '  -1.09%  '
$ws.Range('D48').Value = '85.53'
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('D49').Value = '0.0₇0992'
$ws.Range('E49').Value = '  -1.27%  '
$ws.Range('D50').Value = '0.0955'
$ws.Range('E50').Value = '  -2.04%  '
$ws.Range('D51').Value = '0.0494'
$ws.Range('E51').Value = '  -0.70%  '
